# RPAR_holdings.xlsx update:
#   - bump the "as of" date in the confidential disclosure blurb (A18)
#     from 2021-03-25 to 2021-03-26
#   - refresh the Weight (D) / Percent Change (E) figures for rows 2-15
#
# The sheet ships protected, so it has to be unprotected before any
# cell can be written, and is re-protected afterwards to restore the
# original (protected) state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("password")

# --- Disclosure text: 2021-03-25 -> 2021-03-26 ---------------------------
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) refresh, rows 2-15 -------------------
$values = @(
    @{ Row = 2;  D = 0.05498163670169722;  E = 0.01634260398720144 }
    @{ Row = 3;  D = 0.02315101749181974;  E = 0.01370703764320758 }
    @{ Row = 4;  D = 0.03083735016387586;  E = 0.02488151658767768 }
    @{ Row = 5;  D = 0.03174062957560077;  E = 0.02106115836370992 }
    @{ Row = 6;  D = 0.03246122326360603;  E = 0.04652326163081555 }
    @{ Row = 7;  D = 0.01868942615031111;  E = 0.01533198175367434 }
    @{ Row = 8;  D = 0.004569172934365756; E = 0.0199911150599732 }
    @{ Row = 9;  D = 0.006533382771247704; E = 0.02112676056338025 }
    @{ Row = 10; D = 0.06962355915093088;  E = 0.00291545189504383 }
    @{ Row = 11; D = 0.06978594646090391;  E = 0.002908667830133549 }
    @{ Row = 12; D = 0.1484653045979967;   E = -0.00350007291818577 }
    @{ Row = 13; D = 0.3932452291961558;   E = -0.003334503334503358 }
    @{ Row = 14; D = 0.1159161215414886;   E = -0.0006128953174797669 }
)

foreach ($item in $values) {
    $ws.Range("D$($item.Row)").Value = $item.D
    $ws.Range("E$($item.Row)").Value = $item.E
}

# Row 15 only has its Percent Change (E) figure refreshed.
$ws.Range("E15").Value = 0.003181776354783894

$ws.Protect("password")
